{"js": "// Update the worksheet date and all two-digit \u00f7 one-digit division\n// problems/answers to the new values, per the commit's regenerated\n// content. Each \"from\" text is unique in the document, so a simple\n// matchCase search + full-text replace per pair is safe and avoids\n// any row/column/table-position bookkeeping.\nconst replacements = [\n  [\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"],\n  [\"38\u00f75=7, 3\", \"75\u00f78=9, 3\"],\n  [\"39\u00f72=19, 1\", \"17\u00f79=1, 8\"],\n  [\"79\u00f76=13, 1\", \"75\u00f77=10, 5\"],\n  [\"98\u00f72=49, 0\", \"40\u00f78=5, 0\"],\n  [\"91\u00f75=18, 1\", \"76\u00f79=8, 4\"],\n  [\"15\u00f76=2, 3\", \"61\u00f75=12, 1\"],\n  [\"87\u00f73=29, 0\", \"73\u00f74=18, 1\"],\n  [\"39\u00f77=5, 4\", \"47\u00f79=5, 2\"],\n  [\"14\u00f74=3, 2\", \"57\u00f76=9, 3\"],\n  [\"23\u00f77=3, 2\", \"32\u00f74=8, 0\"],\n  [\"33\u00f79=3, 6\", \"30\u00f79=3, 3\"],\n  [\"94\u00f73=31, 1\", \"51\u00f77=7, 2\"],\n  [\"50\u00f79=5, 5\", \"34\u00f74=8, 2\"],\n  [\"46\u00f75=9, 1\", \"90\u00f75=18, 0\"],\n  [\"80\u00f79=8, 8\", \"51\u00f76=8, 3\"],\n  [\"37\u00f77=5, 2\", \"40\u00f78=5, 0\"],\n  [\"86\u00f77=12, 2\", \"86\u00f72=43, 0\"],\n  [\"69\u00f76=11, 3\", \"17\u00f79=1, 8\"],\n  [\"72\u00f79=8, 0\", \"88\u00f76=14, 4\"],\n  [\"99\u00f76=16, 3\", \"39\u00f79=4, 3\"],\n  [\"95\u00f79=10, 5\", \"69\u00f79=7, 6\"],\n  [\"80\u00f75=16, 0\", \"46\u00f78=5, 6\"],\n  [\"85\u00f76=14, 1\", \"10\u00f73=3, 1\"],\n  [\"19\u00f78=2, 3\", \"26\u00f78=3, 2\"],\n  [\"87\u00f75=17, 2\", \"61\u00f74=15, 1\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all two-digit \u00f7 one-digit division\n# problems/answers to the new values, per the commit's regenerated\n# content. Each \"from\" text is unique in the document, so a fresh\n# Content.Find.Execute replace-all per pair is safe and avoids any\n# row/column/table-position bookkeeping.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-02 Tuesday\", \"2024-04-03 Wednesday\"),\n    @(\"38\u00f75=7, 3\", \"75\u00f78=9, 3\"),\n    @(\"39\u00f72=19, 1\", \"17\u00f79=1, 8\"),\n    @(\"79\u00f76=13, 1\", \"75\u00f77=10, 5\"),\n    @(\"98\u00f72=49, 0\", \"40\u00f78=5, 0\"),\n    @(\"91\u00f75=18, 1\", \"76\u00f79=8, 4\"),\n    @(\"15\u00f76=2, 3\", \"61\u00f75=12, 1\"),\n    @(\"87\u00f73=29, 0\", \"73\u00f74=18, 1\"),\n    @(\"39\u00f77=5, 4\", \"47\u00f79=5, 2\"),\n    @(\"14\u00f74=3, 2\", \"57\u00f76=9, 3\"),\n    @(\"23\u00f77=3, 2\", \"32\u00f74=8, 0\"),\n    @(\"33\u00f79=3, 6\", \"30\u00f79=3, 3\"),\n    @(\"94\u00f73=31, 1\", \"51\u00f77=7, 2\"),\n    @(\"50\u00f79=5, 5\", \"34\u00f74=8, 2\"),\n    @(\"46\u00f75=9, 1\", \"90\u00f75=18, 0\"),\n    @(\"80\u00f79=8, 8\", \"51\u00f76=8, 3\"),\n    @(\"37\u00f77=5, 2\", \"40\u00f78=5, 0\"),\n    @(\"86\u00f77=12, 2\", \"86\u00f72=43, 0\"),\n    @(\"69\u00f76=11, 3\", \"17\u00f79=1, 8\"),\n    @(\"72\u00f79=8, 0\", \"88\u00f76=14, 4\"),\n    @(\"99\u00f76=16, 3\", \"39\u00f79=4, 3\"),\n    @(\"95\u00f79=10, 5\", \"69\u00f79=7, 6\"),\n    @(\"80\u00f75=16, 0\", \"46\u00f78=5, 6\"),\n    @(\"85\u00f76=14, 1\", \"10\u00f73=3, 1\"),\n    @(\"19\u00f78=2, 3\", \"26\u00f78=3, 2\"),\n    @(\"87\u00f75=17, 2\", \"61\u00f74=15, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
